# RDCC-5182 Added Version check
#
# Adds a new "VERSION" worksheet right after the existing
# "Service to CW Roles Mapping" sheet, containing a small
# File version / vx.xx label pair in row 6, and makes that
# new sheet the active tab.

$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after the first (existing) sheet.
$versionSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$versionSheet.Name = "VERSION"

# Populate row 6 with the version info, matching the target layout.
$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Make B6 the selected cell on the new sheet and make it the active tab.
$versionSheet.Range("B6").Select()
